$d = $word.ActiveDocument

function Append-ParaText {
    param($doc, $paraIndex, $text)
    $p = $doc.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $rc = $doc.Range($r.End - 1, $r.End - 1)
    $rc.InsertAfter($text)
}

function Replace-ParaText {
    param($doc, $paraIndex, $oldText, $newText)
    $p = $doc.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# Work from the bottom of the document upward so that paragraph indices
# for not-yet-processed paragraphs stay stable even after the new
# paragraph is inserted near the top (after paragraph 2).

# "Permettre la résolution des noms" -> append " 28 CHARLEMAGNE"
Append-ParaText $d 35 " 28 CHARLEMAGNE"

# "Permettre à tout utilisateur authentifié d’accéder à ses données" -> append " 27 CHARLEMAGNE"
Append-ParaText $d 34 " 27 CHARLEMAGNE"

# "Permettre à tout utilisateur d’accéder au site et aux applicatifs" -> append " 26 CHARLEMAGNE"
Append-ParaText $d 33 " 26 CHARLEMAGNE"

# "Sécuriser l’ensemble de l’infrastructure et du réseau" -> append " 25 CHARLEMAGNE"
Append-ParaText $d 32 " 25 CHARLEMAGNE"

# "Gestion des utilisateurs et des accès" -> append " 24 CHARLEMAGNE"
Append-ParaText $d 30 " 24 CHARLEMAGNE"

# "Suivi des dépôts et validation" -> append " 23 CHARLEMAGNE"
Append-ParaText $d 29 " 23 CHARLEMAGNE"

# "Espace administrateur" -> append " 22 CHARLEMAGNE"
Append-ParaText $d 28 " 22 CHARLEMAGNE"

# "Suivi de validation des défis" -> append " 21 BERENGUER"
Append-ParaText $d 26 " 21 BERENGUER"

# "Gestion des données et coordonnées " -> append "20 BERENGUER"
Append-ParaText $d 25 "20 BERENGUER"

# "Listing des défis souscrits et description" -> append " 19 BERENGUER" (paragraph-scoped to avoid the duplicate at paragraph 5)
Replace-ParaText $d 24 "Listing des défis souscrits et description" "Listing des défis souscrits et description 19 BERENGUER"

# "Déposer des vidéos de défis" -> append " 18 BERENGUER" (paragraph-scoped to avoid the duplicate at paragraph 4)
Append-ParaText $d 23 " 18 BERENGUER"

# "Espace Client/utilisateur " -> append "17 BERENGUER"
Append-ParaText $d 22 "17 BERENGUER"

# "Connexion " -> append "16 BERENGUER"
Append-ParaText $d 20 "16 BERENGUER"

# "Paiement en ligne" -> append " 15 BERENGUER"
Append-ParaText $d 19 " 15 BERENGUER"

# "Formulaire d’inscription " -> append " 14 BERENGUER"
Append-ParaText $d 18 " 14 BERENGUER"

# "Défi démo" -> append " 13 BERENGUER"
Append-ParaText $d 17 " 13 BERENGUER"

# "Classement général" -> append " 12 BERENGUER"
Append-ParaText $d 16 " 12 BERENGUER"

# "Présentation du site" -> append " 11 BERENGUER"
Append-ParaText $d 15 " 11 BERENGUER"

# "Accueil" -> append " 10 BERENGUER"
Append-ParaText $d 14 " 10 BERENGUER"

# "Espace public" -> append " 09 BERENGUER"
Append-ParaText $d 13 " 09 BERENGUER"

# "Back-end" -> append " 08 BERENGUER"
Append-ParaText $d 11 " 08 BERENGUER"

# "Espace de jeu" -> append " 07 BERENGUER"
Append-ParaText $d 10 " 07 BERENGUER"

# "Front-end" -> append " 06 BERENGUER"
Append-ParaText $d 9 " 06 BERENGUER"

# "Suivi de validation des défis" -> append " 05 CHARLEMAGNE"
Append-ParaText $d 7 " 05 CHARLEMAGNE"

# "Gestion des données et coordonnées" -> append " 04 CHARLEMAGNE"
Append-ParaText $d 6 " 04 CHARLEMAGNE"

# "Listing des défis souscrits et description" -> append " 03 CHARLEMAGNE"
Replace-ParaText $d 5 "Listing des défis souscrits et description" "Listing des défis souscrits et description 03 CHARLEMAGNE"

# "Déposer des vidéos de défis" -> append " 02 CHARLEMAGNE"
Append-ParaText $d 4 " 02 CHARLEMAGNE"

# "Documentation " -> "Documentation initiale 00 CHARLEMAGNE" then insert a new
# paragraph below it: "Documentation suivie 01 CHARLEMAGNE"
Replace-ParaText $d 2 "Documentation " "Documentation initiale 00 CHARLEMAGNE"

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$rc = $d.Range($p3.Range.Start, $p3.Range.Start)
$rc.InsertAfter("Documentation suivie 01 CHARLEMAGNE")
